$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) text changes.
#    Old: A title | B title_short | C navbar_logo | D primary_colour |
#         E secondary_colour | F favicon_link | G data_schemas | H partner_name |
#         I partner_website | J google_analytics_key | K "G-tag script"
#    New: A title | B title_short | C favicon_link | D data_schemas |
#         E partner_name | F partner_website | G google_analytics_key |
#         H gTag_script | I organization_name | J organization_link |
#         K logo_colour | L footer_colour | M body_foreground_colour |
#         N body_background_colour | O find_candidates_button | P candidate_names_colour
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "favicon_link"
$ws.Range("D1").Value = "data_schemas"
$ws.Range("E1").Value = "partner_name"
$ws.Range("F1").Value = "partner_website"
$ws.Range("G1").Value = "google_analytics_key"
$ws.Range("H1").Value = "gTag_script"
$ws.Range("I1").Value = "organization_name"
$ws.Range("J1").Value = "organization_link"
$ws.Range("K1").Value = "logo_colour"
$ws.Range("L1").Value = "footer_colour"
$ws.Range("M1").Value = "body_foreground_colour"
$ws.Range("N1").Value = "body_background_colour"
$ws.Range("O1").Value = "find_candidates_button"
$ws.Range("P1").Value = "candidate_names_colour"

# Remove the old trailing style-only placeholder cells (L1:Y1 in the old sheet);
# only A1:P1 should remain populated on row 1.
$ws.Range("Q1:Y1").Clear()

# ---------------------------------------------------------------------------
# 2. New data values for the newly added colour columns on row 2.
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = "#FFC4AB"
$ws.Range("L2").Value = "#1D3437"
$ws.Range("M2").Value = "#30474A"
$ws.Range("N2").Value = "#E5E5E5"
$ws.Range("O2").Value = "#E07D54"
$ws.Range("P2").Value = "#F2AA71"

Write-Output "values done"
